$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing "C2" column (H) for a new "C1" column,
# then clear the auto-styled blank data cells so only the header remains.
$ws.Columns("H").Insert()
$ws.Range("H2:H4").Clear()
$ws.Range("H1").Value = "C1"

# Insert a new column before the (now shifted) "D2" column (K) for a new "D1"
# column, then clear the auto-styled blank data cells so only the header remains.
$ws.Columns("K").Insert()
$ws.Range("K2:K4").Clear()
$ws.Range("K1").Value = "D1"

# Update the active selection to match the saved view state.
$null = $ws.Range("E6").Select()
